# Insert 3 new weekly-price rows for "Comercializadora del Agro de Limarí - Frutilla"
# ahead of the existing rows 304-361 (shifting that whole historic block down by
# 3 rows, to 307-364), then populate the newly inserted rows 304:306 with the
# latest week's observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (rows 304:361) down by three rows.
$ws.Rows("304:306").Insert()

# Row 304 - Especial
$ws.Range("A304").Value = 2
$ws.Range("B304").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C304").Value = "Coquimbo"
$ws.Range("D304").Value = 44637
$ws.Range("E304").Value = 4
$ws.Range("F304").Value = "Fruta"
$ws.Range("G304").Value = 100101
$ws.Range("H304").Value = "Berries"
$ws.Range("I304").Value = 100112025
$ws.Range("J304").Value = "Frutilla"
$ws.Range("K304").Value = "Sin especificar"
$ws.Range("L304").Value = "Especial"
$ws.Range("M304").Value = 400
$ws.Range("N304").Value = 12000
$ws.Range("O304").Value = 12500
$ws.Range("P304").Value = 12250
$ws.Range("Q304").Value = "$/bandeja 7 kilos"
$ws.Range("R304").Value = "Provincia de Melipilla"
$ws.Range("S304").Value = 1750
$ws.Range("T304").Value = 7

# Row 305 - Primera
$ws.Range("A305").Value = 2
$ws.Range("B305").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C305").Value = "Coquimbo"
$ws.Range("D305").Value = 44637
$ws.Range("E305").Value = 4
$ws.Range("F305").Value = "Fruta"
$ws.Range("G305").Value = 100101
$ws.Range("H305").Value = "Berries"
$ws.Range("I305").Value = 100112025
$ws.Range("J305").Value = "Frutilla"
$ws.Range("K305").Value = "Sin especificar"
$ws.Range("L305").Value = "Primera"
$ws.Range("M305").Value = 500
$ws.Range("N305").Value = 10000
$ws.Range("O305").Value = 10500
$ws.Range("P305").Value = 10250
$ws.Range("Q305").Value = "$/bandeja 7 kilos"
$ws.Range("R305").Value = "Provincia de Melipilla"
$ws.Range("S305").Value = 1464
$ws.Range("T305").Value = 7

# Row 306 - Segunda
$ws.Range("A306").Value = 2
$ws.Range("B306").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C306").Value = "Coquimbo"
$ws.Range("D306").Value = 44637
$ws.Range("E306").Value = 4
$ws.Range("F306").Value = "Fruta"
$ws.Range("G306").Value = 100101
$ws.Range("H306").Value = "Berries"
$ws.Range("I306").Value = 100112025
$ws.Range("J306").Value = "Frutilla"
$ws.Range("K306").Value = "Sin especificar"
$ws.Range("L306").Value = "Segunda"
$ws.Range("M306").Value = 400
$ws.Range("N306").Value = 8000
$ws.Range("O306").Value = 8500
$ws.Range("P306").Value = 8250
$ws.Range("Q306").Value = "$/bandeja 7 kilos"
$ws.Range("R306").Value = "Provincia de Melipilla"
$ws.Range("S306").Value = 1179
$ws.Range("T306").Value = 7
